$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "51.083.28"
$ws.Range("E2").Value = "  -1.03%  "
$ws.Range("D3").Value = "2.942.59"
$ws.Range("E3").Value = "  -1.43%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.00"
$ws.Range("E4").Value = "  +0.12%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "374.38"
$ws.Range("E5").Value = "  -1.45%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "101.23"
$ws.Range("E7").Value = "  -1.55%  "
$ws.Range("E8").Value = "  +0.08%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.584"
$ws.Range("E9").Value = "  -1.84%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "36.28"
$ws.Range("E10").Value = "  -2.61%  "
$ws.Range("E11").Value = "  -0.78%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.0849"
$ws.Range("E12").Value = "  +0.32%  "
$ws.Range("D13").Value = "3.408.53"
$ws.Range("E13").Value = "  -1.26%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "18.03"
$ws.Range("E14").Value = "  -2.21%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "7.53"
$ws.Range("E15").Value = "  -0.78%  "
$ws.Range("D16").Value = "2.933.45"
$ws.Range("E16").Value = "  -1.77%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "11.03"
$ws.Range("E17").Value = "  +48.39%  "
$ws.Range("E18").Value = "  -0.07%  "
$ws.Range("D19").Value = "51.065.03"
$ws.Range("E19").Value = "  -0.98%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "3.14"
$ws.Range("E20").Value = "  -5.65%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "12.42"
$ws.Range("E21").Value = "  -4.32%  "
$ws.Range("E22").Value = "  -0.67%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "265.04"
$ws.Range("E23").Value = "  +1.03%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "68.68"
$ws.Range("E24").Value = "  -0.94%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "3.13"
$ws.Range("E25").Value = "  +9.66%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "8.18"
$ws.Range("E26").Value = "  -0.56%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "7.56"
$ws.Range("E27").Value = "  -2.54%  "
$ws.Range("E28").Value = "  -0.03%  "
$ws.Range("B29").Value = "EthereumClassic"
$ws.Range("C29").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "25.64"
$ws.Range("E29").Value = "  -1.13%  "
$ws.Range("B30").Value = "Kaspa"
$ws.Range("C30").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.164"
$ws.Range("E30").Value = "  -4.18%  "
$ws.Range("B31").Value = "Hedera"
$ws.Range("C31").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.110"
$ws.Range("E31").Value = "  -5.75%  "
$ws.Range("E32").Value = "  +1.26%  "
$ws.Range("E33").Value = "  -0.82%  "
$ws.Range("E34").Value = "  -1.83%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "33.43"
$ws.Range("E35").Value = "  -4.79%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.0444"
$ws.Range("E36").Value = "  -0.57%  "
$ws.Range("E37").Value = "  -0.13%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "3.18"
$ws.Range("E38").Value = "  +4.28%  "
$ws.Range("E39").Value = "  -0.41%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "16.33"
$ws.Range("E40").Value = "  -4.83%  "
$ws.Range("E41").Value = "  -3.77%  "
$ws.Range("E42").Value = "  -3.39%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "119.85"
$ws.Range("E43").Value = "  -4.53%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "21.39"
$ws.Range("E44").Value = "  -1.60%  "
$ws.Range("B45").Value = "TheGraph"
$ws.Range("C45").Value = "https://coinranking.com/coin/qhd1biQ7M+thegraph-grt"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.274"
$ws.Range("E45").Value = "  -3.45%  "
$ws.Range("B46").Value = "WEMIXToken"
$ws.Range("C46").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.03"
$ws.Range("E46").Value = "  -0.96%  "
$ws.Range("B47").Value = "NEARProtocol"
$ws.Range("C47").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "3.34"
$ws.Range("E47").Value = "  +2.61%  "
$ws.Range("B48").Value = "ApeXProtocol"
$ws.Range("C48").Value = "https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "2.30"
$ws.Range("E48").Value = "  -2.96%  "
$ws.Range("B49").Value = "Maker"
$ws.Range("C49").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D49").Value = "1.996.89"
$ws.Range("E49").Value = "  -2.00%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.0325"
$ws.Range("E50").Value = "  -2.99%  "
$ws.Range("E51").Value = "  +1.22%  "
